$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.905.12'
$ws.Range("E2").Value = '  -0.33%  '

# Row 3
$ws.Range("D3").Value = '3.763.40'
$ws.Range("E3").Value = '  -1.41%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '637.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '

# Row 7
$ws.Range("D7").Value = '3.762.32'
$ws.Range("E7").Value = '  -1.37%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("E9").Value = '  +0.17%  '

# Row 10
$ws.Range("E10").Value = '  -2.56%  '

# Row 11
$ws.Range("E11").Value = '  +0.42%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.69%  '

# Row 13
$ws.Range("E13").Value = '  -4.89%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.29%  '

# Row 15
$ws.Range("D15").Value = '4.395.00'
$ws.Range("E15").Value = '  -1.40%  '

# Row 16
$ws.Range("D16").Value = '3.759.22'
$ws.Range("E16").Value = '  -2.10%  '

# Row 17
$ws.Range("D17").Value = '68.827.33'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.83%  '

# Row 19
$ws.Range("E19").Value = '  -0.16%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '470.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.85%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.03%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.703'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.86%  '

# Row 24
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000143'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.46%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.41%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.27%  '

# Row 27
$ws.Range("E27").Value = '  -2.45%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("D30").Value = '3.910.53'
$ws.Range("E30").Value = '  -1.43%  '

# Row 31
$ws.Range("E31").Value = '  -0.80%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.57%  '

# Row 33
$ws.Range("E33").Value = '  -2.82%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.74%  '

# Row 35
$ws.Range("E35").Value = '  +16.60%  '

# Row 36
$ws.Range("E36").Value = '  +0.10%  '

# Row 37
$ws.Range("D37").Value = '3.717.96'
$ws.Range("E37").Value = '  -1.16%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.25%  '

# Row 39
$ws.Range("E39").Value = '  -1.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.75%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.66%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.15%  '

# Row 43
$ws.Range("E43").Value = '  -2.65%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '44.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.84%  '

# Row 46
$ws.Range("E46").Value = '  +3.75%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '155.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.53%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.95%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.25%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.293'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.38%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.22%  '
